# Update the "想去人数" (interest count) figures in column F across sheets,
# reflecting refreshed stats from the data source.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 16521
$ws1.Range("F4").Value  = 22
$ws1.Range("F5").Value  = 747
$ws1.Range("F6").Value  = 15639
$ws1.Range("F7").Value  = 79
$ws1.Range("F10").Value = 110
$ws1.Range("F11").Value = 658
$ws1.Range("F14").Value = 86
$ws1.Range("F15").Value = 1171
$ws1.Range("F16").Value = 18
$ws1.Range("F17").Value = 29
$ws1.Range("F18").Value = 39
$ws1.Range("F19").Value = 564
$ws1.Range("F20").Value = 47
$ws1.Range("F21").Value = 49
$ws1.Range("F24").Value = 5
$ws1.Range("F25").Value = 78
$ws1.Range("F26").Value = 286
$ws1.Range("F27").Value = 387
$ws1.Range("F30").Value = 5865
$ws1.Range("F31").Value = 5274

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 86

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 16521
$ws4.Range("F4").Value  = 22
$ws4.Range("F5").Value  = 747
$ws4.Range("F6").Value  = 15639
$ws4.Range("F7").Value  = 79
$ws4.Range("F10").Value = 110
$ws4.Range("F11").Value = 658
$ws4.Range("F14").Value = 86
$ws4.Range("F15").Value = 1171
$ws4.Range("F16").Value = 18
$ws4.Range("F18").Value = 39
$ws4.Range("F19").Value = 564
$ws4.Range("F20").Value = 47
$ws4.Range("F21").Value = 49
$ws4.Range("F22").Value = 86
$ws4.Range("F26").Value = 5
$ws4.Range("F27").Value = 78
$ws4.Range("F28").Value = 286
$ws4.Range("F29").Value = 387
$ws4.Range("F32").Value = 5865
$ws4.Range("F34").Value = 5274
